$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This change ("Generate Report for Handoff") inserts a new handoff entry for
# file 0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md on every sheet, just above the
# existing ".localization-config" row (which shifts down by one row).
# ---------------------------------------------------------------------------

$newGuid   = "0bcd2ce7-5fcc-48b5-94a0-47d35f7452df"
$newMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/68b703065576b0821fc5d6d8974bfd03bfe02968/e2e/$newGuid.md"

# ============================= Sheet "Overview" =============================
$ws = $wb.Worksheets.Item("Overview")

# Move the ".localization-config" row (old row 4) down to row 5, preserving
# its values and hyperlink.
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("C5").Value = "Not to be localized"

# Overwrite row 4 with the new handoff entry.
$ws.Range("A4").Value = "$newGuid.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"

# Fix up hyperlinks: A4 now points at the new .md file, A5 keeps the old
# ".localization-config" link.
$ws.Range("A4").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A4"), $newMdUrl, "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/68b703065576b0821fc5d6d8974bfd03bfe02968/.localization-config", "", "", ".localization-config")

# ============================= Sheet "zh-cn" ================================
$ws = $wb.Worksheets.Item("zh-cn")

$zhXlf    = "$newGuid.a1ce4849b9565447c4a185bf00b08dbaf7048def.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/abb3b05685450e9f50bf2a62cc43313e624cf6f4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf"

# Move the ".localization-config" row (old row 4) down to row 5.
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# Overwrite row 4 with the new handoff entry.
$ws.Range("A4").Value = "$newGuid.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = $zhXlf
$ws.Range("D4").Value = "2016-03-02 09:25:17"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Include"

# Fix up hyperlinks.
$ws.Range("A4").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A4"), $newMdUrl, "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("C4"), $zhXlfUrl, "", "", $zhXlf)
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/68b703065576b0821fc5d6d8974bfd03bfe02968/.localization-config", "", "", ".localization-config")

# ============================= Sheet "de-de" ================================
$ws = $wb.Worksheets.Item("de-de")

$deXlf    = "$newGuid.a1ce4849b9565447c4a185bf00b08dbaf7048def.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd8ddf4325ab0f5e12e25a4e46a9e1fd82e30a8e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf"

# Move the ".localization-config" row (old row 4) down to row 5.
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# Overwrite row 4 with the new handoff entry.
$ws.Range("A4").Value = "$newGuid.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = $deXlf
$ws.Range("D4").Value = "2016-03-02 09:25:29"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Include"

# Fix up hyperlinks.
$ws.Range("A4").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A4"), $newMdUrl, "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("C4"), $deXlfUrl, "", "", $deXlf)
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/68b703065576b0821fc5d6d8974bfd03bfe02968/.localization-config", "", "", ".localization-config")
